# Applies the edits described by the commit:
# "cabezales de las actividades arreglados"
#
# 1) Merge the "IEEE 1074" + " " runs into a single run "IEEE 1074 ".
# 2) Fix the delivery date from 25/6/2019 to 26/6/2019.
# 3) Move the "_GoBack" bookmark from the "Testeo" section's
#    "SIN ACTIVIDADES POR EL MOMENTO" paragraph to the start of the
#    "Ruta en GitLab" paragraph near the top of the document.
# 4) Merge the three runs (plus removing the gramatical-error proof
#    marks) describing the testing activities into a single run.

$d = $word.ActiveDocument

# --- 1) Merge "IEEE 1074" + " " into a single run -------------------------
$d.Content.Find.Execute("IEEE 1074 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "IEEE 1074 ", 2) | Out-Null

# --- 2) Fix the date -------------------------------------------------------
$d.Content.Find.Execute("Primera entrega 25/6/2019", $true, $false, $false, `
                         $false, $false, $true, 1, $false, `
                         "Primera entrega 26/6/2019", 2) | Out-Null

# --- 3) Move the _GoBack bookmark ------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$target = $d.Paragraphs(3)
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)
$d.Bookmarks.Add("_GoBack", $insertionPoint) | Out-Null

# --- 4) Merge the "realizaran" runs into a single run -----------------------
$testText = "En esta etapa se documentarán las actividades que se realizaran" + `
    " con respecto al testeo del software previo a la implementación del " + `
    "mismo en las terminales del cliente con el objetivo de identificar " + `
    "errores e incongruencias"
$d.Content.Find.Execute($testText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $testText, 2) | Out-Null
